$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.643745422363281
$ws.Range("B1").Value = 1.765841841697693
$ws.Range("C1").Value = 2.010205984115601
$ws.Range("D1").Value = 3.206500291824341
$ws.Range("E1").Value = 3.621964454650879
